$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 currently holds the phone number "09876543" stored as text in
# column A (with 0 points in column C). The fix re-records that phone
# number as a genuine number (9876543, dropping the leading zero) and
# keeps the original text record as history by pushing it down to a new
# row 44 with the same B/C contents it always had.
#
# Inserting a whole row at 43 shifts the existing row 43 (untouched,
# verbatim) down to become row 44 - that's the cleanest way to get an
# exact duplicate of the old row's cells (including the blank-looking
# birthday cell) without retyping them. Only then do we overwrite the
# now-empty row 43 with the corrected numeric phone + its points value.
$ws.Rows("43:43").Insert()

$ws.Range("A43").Value = 9876543
$ws.Range("C43").Value = 0
